$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE BALANCE")
$ws.Range("B76").Value = "UT(0-0-1)"
$ws.Range("G76").Formula = '=IF(ISBLANK(Table13[[#This Row],[EARNED]]),"",Table13[[#This Row],[EARNED]])'
Write-Host "done"
